$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 2 updates
$ws.Range("F2").Value = "DSI-BUILD-OAR-50"
$ws.Range("I2").Value = "HW - Screen for PC"
$ws.Range("J2").Value = "S00001610001"
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 5000
$ws.Range("U2").Value = "RENT FIXED COST"
$ws.Range("V2").Value = "S"
$ws.Range("W2").Value = 99

# Row 3 updates
$ws.Range("E3").Value = 362
$ws.Range("F3").Value = "DSI-BUILD-OAR-51"
$ws.Range("I3").Value = "SW & Solutions - IS-HR"
$ws.Range("J3").Value = "S00001610001"
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10000
$ws.Range("U3").Value = "SAAS FIXED COST"
$ws.Range("V3").Value = "S"
$ws.Range("W3").Value = 99
$ws.Range("Y3").Value = "yes"
$ws.Range("AB3").Value = "EMEAAD\lbaisin"
$ws.Range("AE3").Value = "EMEAAD\srofidal"

# Selection update: active cell moves to E3, no frozen/top-left scroll cell
$ws.Range("E3").Select()
